# Updates the account-statement ("Estado de Cuenta") period rows.
# Row 16 and Row 18 swap their "Periodo Mora" (E) and "Valor Mora" (F) data:
#   Row16: 1802/42530  ->  1809/42955
#   Row17: 1808/42955  ->  unchanged
#   Row18: 1809/42955  ->  1802/42530

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("E16").Value = "1809"
$ws.Range("F16").Value = 42955

$ws.Range("E18").Value = "1802"
$ws.Range("F18").Value = 42530

$wb.Save()
